# Insert two new data rows at the top of the data block (row 85), pushing the
# existing rows 85-132 down to 87-134, then populate the two new rows with
# their data (a new weekly "Mora" price report entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 85 (shifts rows 85:132 -> 87:134)
$ws.Range("A85:A86").EntireRow.Insert()

# --- New row 85 ---
$ws.Cells.Item(85, 1).Value  = 6
$ws.Cells.Item(85, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(85, 3).Value  = "Metropolitana"
$ws.Cells.Item(85, 4).Value  = 44957
$ws.Cells.Item(85, 5).Value  = 13
$ws.Cells.Item(85, 6).Value  = "Fruta"
$ws.Cells.Item(85, 7).Value  = 100101
$ws.Cells.Item(85, 8).Value  = "Berries"
$ws.Cells.Item(85, 9).Value  = 100101008
$ws.Cells.Item(85, 10).Value = "Mora"
$ws.Cells.Item(85, 11).Value = "Sin especificar"
$ws.Cells.Item(85, 12).Value = "Primera"
$ws.Cells.Item(85, 13).Value = 250
$ws.Cells.Item(85, 14).Value = 3000
$ws.Cells.Item(85, 15).Value = 3000
$ws.Cells.Item(85, 16).Value = 3000
$ws.Cells.Item(85, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(85, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(85, 19).Value = 1500
$ws.Cells.Item(85, 20).Value = 2

# --- New row 86 ---
$ws.Cells.Item(86, 1).Value  = 6
$ws.Cells.Item(86, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(86, 3).Value  = "Metropolitana"
$ws.Cells.Item(86, 4).Value  = 44957
$ws.Cells.Item(86, 5).Value  = 13
$ws.Cells.Item(86, 6).Value  = "Fruta"
$ws.Cells.Item(86, 7).Value  = 100101
$ws.Cells.Item(86, 8).Value  = "Berries"
$ws.Cells.Item(86, 9).Value  = 100101008
$ws.Cells.Item(86, 10).Value = "Mora"
$ws.Cells.Item(86, 11).Value = "Sin especificar"
$ws.Cells.Item(86, 12).Value = "Primera"
$ws.Cells.Item(86, 13).Value = 300
$ws.Cells.Item(86, 14).Value = 3000
$ws.Cells.Item(86, 15).Value = 3000
$ws.Cells.Item(86, 16).Value = 3000
$ws.Cells.Item(86, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(86, 18).Value = "Región del Maule"
$ws.Cells.Item(86, 19).Value = 1500
$ws.Cells.Item(86, 20).Value = 2
